# The workbook gains one new data row: a new weekly price observation is
# inserted as row 3 (pushing the existing rows 3-21 down to rows 4-22,
# which keeps all of their original data intact).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts rows 3..21 down to 4..22.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new observation's data.
$ws.Cells.Item(3, 1).Value  = 1
$ws.Cells.Item(3, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(3, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(3, 4).Value  = 44963
$ws.Cells.Item(3, 5).Value  = 15
$ws.Cells.Item(3, 6).Value  = 100112001
$ws.Cells.Item(3, 7).Value  = 'Berenjena'
$ws.Cells.Item(3, 8).Value  = 'Sin especificar'
$ws.Cells.Item(3, 9).Value  = 'Primera'
$ws.Cells.Item(3, 10).Value = 130
$ws.Cells.Item(3, 11).Value = 4000
$ws.Cells.Item(3, 12).Value = 4500
$ws.Cells.Item(3, 13).Value = 4250
$ws.Cells.Item(3, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(3, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(3, 16).Value = 71
$ws.Cells.Item(3, 17).Value = 60
$ws.Cells.Item(3, 18).Value = 'Hortaliza'
